# Apply updated NATMI LR-pair statistics (Dcn-Met) per Dr Hou's advice
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 22.73180333333334
$ws.Cells.Item(2, 8).Value = 68.19541000000001
$ws.Cells.Item(2, 9).Value = 0.007290521456144306
$ws.Cells.Item(2, 10).Value = 0.007290521456144306
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.167735333333333
$ws.Cells.Item(2, 14).Value = 9.503206
$ws.Cells.Item(2, 15).Value = 0.1182666224938439
$ws.Cells.Item(2, 16).Value = 0.1182666224938439
$ws.Cells.Item(2, 17).Value = 72.00833660938446
$ws.Cells.Item(2, 18).Value = 648.0750294844602
$ws.Cells.Item(2, 19).Value = 0.0008622253488370877
$ws.Cells.Item(2, 20).Value = 0.0008622253488370878

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 22.73180333333334
$ws.Cells.Item(3, 8).Value = 68.19541000000001
$ws.Cells.Item(3, 9).Value = 0.007290521456144306
$ws.Cells.Item(3, 10).Value = 0.007290521456144306
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.9421210000000001
$ws.Cells.Item(3, 14).Value = 2.826363
$ws.Cells.Item(3, 15).Value = 0.03517385669126484
$ws.Cells.Item(3, 16).Value = 0.03517385669126484
$ws.Cells.Item(3, 17).Value = 21.41610928820334
$ws.Cells.Item(3, 18).Value = 192.74498359383
$ws.Cells.Item(3, 19).Value = 0.0002564357569030113
$ws.Cells.Item(3, 20).Value = 0.0002564357569030113

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 22.73180333333334
$ws.Cells.Item(4, 8).Value = 68.19541000000001
$ws.Cells.Item(4, 9).Value = 0.007290521456144306
$ws.Cells.Item(4, 10).Value = 0.007290521456144306
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.442127
$ws.Cells.Item(4, 14).Value = 4.326381
$ws.Cells.Item(4, 15).Value = 0.05384145818700961
$ws.Cells.Item(4, 16).Value = 0.0538414581870096
$ws.Cells.Item(4, 17).Value = 32.78214734569
$ws.Cells.Item(4, 18).Value = 295.03932611121
$ws.Cells.Item(4, 19).Value = 0.0003925323061424901
$ws.Cells.Item(4, 20).Value = 0.00039253230614249

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 22.73180333333334
$ws.Cells.Item(5, 8).Value = 68.19541000000001
$ws.Cells.Item(5, 9).Value = 0.007290521456144306
$ws.Cells.Item(5, 10).Value = 0.007290521456144306
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 21.232711
$ws.Cells.Item(5, 14).Value = 63.69813300000001
$ws.Cells.Item(5, 15).Value = 0.7927180626278817
$ws.Cells.Item(5, 16).Value = 0.7927180626278817
$ws.Cells.Item(5, 17).Value = 482.6578106855034
$ws.Cells.Item(5, 18).Value = 4343.920296169531
$ws.Cells.Item(5, 19).Value = 0.005779328044261717
$ws.Cells.Item(5, 20).Value = 0.005779328044261717

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 3042.696044666667
$ws.Cells.Item(6, 8).Value = 9128.088134
$ws.Cells.Item(6, 9).Value = 0.9758504625824999
$ws.Cells.Item(6, 10).Value = 0.9758504625824997
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 3.167735333333333
$ws.Cells.Item(6, 14).Value = 9.503206
$ws.Cells.Item(6, 15).Value = 0.1182666224938439
$ws.Cells.Item(6, 16).Value = 0.1182666224938439
$ws.Cells.Item(6, 17).Value = 9638.455769284179
$ws.Cells.Item(6, 18).Value = 86746.10192355761
$ws.Cells.Item(6, 19).Value = 0.1154105382686875
$ws.Cells.Item(6, 20).Value = 0.1154105382686875

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 3042.696044666667
$ws.Cells.Item(7, 8).Value = 9128.088134
$ws.Cells.Item(7, 9).Value = 0.9758504625824999
$ws.Cells.Item(7, 10).Value = 0.9758504625824997
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.9421210000000001
$ws.Cells.Item(7, 14).Value = 2.826363
$ws.Cells.Item(7, 15).Value = 0.03517385669126484
$ws.Cells.Item(7, 16).Value = 0.03517385669126484
$ws.Cells.Item(7, 17).Value = 2866.587840297405
$ws.Cells.Item(7, 18).Value = 25799.29056267664
$ws.Cells.Item(7, 19).Value = 0.03432442432298135
$ws.Cells.Item(7, 20).Value = 0.03432442432298135

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 3042.696044666667
$ws.Cells.Item(8, 8).Value = 9128.088134
$ws.Cells.Item(8, 9).Value = 0.9758504625824999
$ws.Cells.Item(8, 10).Value = 0.9758504625824997
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.442127
$ws.Cells.Item(8, 14).Value = 4.326381
$ws.Cells.Item(8, 15).Value = 0.05384145818700961
$ws.Cells.Item(8, 16).Value = 0.0538414581870096
$ws.Cells.Item(8, 17).Value = 4387.954118807006
$ws.Cells.Item(8, 18).Value = 39491.58706926305
$ws.Cells.Item(8, 19).Value = 0.05254121187790965
$ws.Cells.Item(8, 20).Value = 0.05254121187790964

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 3042.696044666667
$ws.Cells.Item(9, 8).Value = 9128.088134
$ws.Cells.Item(9, 9).Value = 0.9758504625824999
$ws.Cells.Item(9, 10).Value = 0.9758504625824997
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 21.232711
$ws.Cells.Item(9, 14).Value = 63.69813300000001
$ws.Cells.Item(9, 15).Value = 0.7927180626278817
$ws.Cells.Item(9, 16).Value = 0.7927180626278817
$ws.Cells.Item(9, 17).Value = 64604.68577725043
$ws.Cells.Item(9, 18).Value = 581442.1719952539
$ws.Cells.Item(9, 19).Value = 0.7735742881129214
$ws.Cells.Item(9, 20).Value = 0.7735742881129214

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.93578
$ws.Cells.Item(10, 8).Value = 5.80734
$ws.Cells.Item(10, 9).Value = 0.0006208414447999517
$ws.Cells.Item(10, 10).Value = 0.0006208414447999516
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 3.167735333333333
$ws.Cells.Item(10, 14).Value = 9.503206
$ws.Cells.Item(10, 15).Value = 0.1182666224938439
$ws.Cells.Item(10, 16).Value = 0.1182666224938439
$ws.Cells.Item(10, 17).Value = 6.13203870356
$ws.Cells.Item(10, 18).Value = 55.18834833204
$ws.Cells.Item(10, 19).Value = [double]"7.342482078068851E-05"
$ws.Cells.Item(10, 20).Value = [double]"7.34248207806885E-05"

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.93578
$ws.Cells.Item(11, 8).Value = 5.80734
$ws.Cells.Item(11, 9).Value = 0.0006208414447999517
$ws.Cells.Item(11, 10).Value = 0.0006208414447999516
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.9421210000000001
$ws.Cells.Item(11, 14).Value = 2.826363
$ws.Cells.Item(11, 15).Value = 0.03517385669126484
$ws.Cells.Item(11, 16).Value = 0.03517385669126484
$ws.Cells.Item(11, 17).Value = 1.82373898938
$ws.Cells.Item(11, 18).Value = 16.41365090442
$ws.Cells.Item(11, 19).Value = [double]"2.183738800739131E-05"
$ws.Cells.Item(11, 20).Value = [double]"2.183738800739131E-05"

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 1.93578
$ws.Cells.Item(12, 8).Value = 5.80734
$ws.Cells.Item(12, 9).Value = 0.0006208414447999517
$ws.Cells.Item(12, 10).Value = 0.0006208414447999516
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.442127
$ws.Cells.Item(12, 14).Value = 4.326381
$ws.Cells.Item(12, 15).Value = 0.05384145818700961
$ws.Cells.Item(12, 16).Value = 0.0538414581870096
$ws.Cells.Item(12, 17).Value = 2.79164060406
$ws.Cells.Item(12, 18).Value = 25.12476543654
$ws.Cells.Item(12, 19).Value = [double]"3.342700869095923E-05"
$ws.Cells.Item(12, 20).Value = [double]"3.342700869095922E-05"

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 1.93578
$ws.Cells.Item(13, 8).Value = 5.80734
$ws.Cells.Item(13, 9).Value = 0.0006208414447999517
$ws.Cells.Item(13, 10).Value = 0.0006208414447999516
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 21.232711
$ws.Cells.Item(13, 14).Value = 63.69813300000001
$ws.Cells.Item(13, 15).Value = 0.7927180626278817
$ws.Cells.Item(13, 16).Value = 0.7927180626278817
$ws.Cells.Item(13, 17).Value = 41.10185729958
$ws.Cells.Item(13, 18).Value = 369.91671569622
$ws.Cells.Item(13, 19).Value = 0.0004921522273209127
$ws.Cells.Item(13, 20).Value = 0.0004921522273209126

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 50.63053333333333
$ws.Cells.Item(14, 8).Value = 151.8916
$ws.Cells.Item(14, 9).Value = 0.016238174516556
$ws.Cells.Item(14, 10).Value = 0.016238174516556
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 3.167735333333333
$ws.Cells.Item(14, 14).Value = 9.503206
$ws.Cells.Item(14, 15).Value = 0.1182666224938439
$ws.Cells.Item(14, 16).Value = 0.1182666224938439
$ws.Cells.Item(14, 17).Value = 160.3841293855111
$ws.Cells.Item(14, 18).Value = 1443.4571644696
$ws.Cells.Item(14, 19).Value = 0.001920434055538685
$ws.Cells.Item(14, 20).Value = 0.001920434055538685

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 50.63053333333333
$ws.Cells.Item(15, 8).Value = 151.8916
$ws.Cells.Item(15, 9).Value = 0.016238174516556
$ws.Cells.Item(15, 10).Value = 0.016238174516556
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.9421210000000001
$ws.Cells.Item(15, 14).Value = 2.826363
$ws.Cells.Item(15, 15).Value = 0.03517385669126484
$ws.Cells.Item(15, 16).Value = 0.03517385669126484
$ws.Cells.Item(15, 17).Value = 47.70008869453333
$ws.Cells.Item(15, 18).Value = 429.3007982508
$ws.Cells.Item(15, 19).Value = 0.0005711592233730895
$ws.Cells.Item(15, 20).Value = 0.0005711592233730895

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 50.63053333333333
$ws.Cells.Item(16, 8).Value = 151.8916
$ws.Cells.Item(16, 9).Value = 0.016238174516556
$ws.Cells.Item(16, 10).Value = 0.016238174516556
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 1.442127
$ws.Cells.Item(16, 14).Value = 4.326381
$ws.Cells.Item(16, 15).Value = 0.05384145818700961
$ws.Cells.Item(16, 16).Value = 0.0538414581870096
$ws.Cells.Item(16, 17).Value = 73.01565914439999
$ws.Cells.Item(16, 18).Value = 657.1409322995999
$ws.Cells.Item(16, 19).Value = 0.0008742869942665149
$ws.Cells.Item(16, 20).Value = 0.0008742869942665148

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 50.63053333333333
$ws.Cells.Item(17, 8).Value = 151.8916
$ws.Cells.Item(17, 9).Value = 0.016238174516556
$ws.Cells.Item(17, 10).Value = 0.016238174516556
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 21.232711
$ws.Cells.Item(17, 14).Value = 63.69813300000001
$ws.Cells.Item(17, 15).Value = 0.7927180626278817
$ws.Cells.Item(17, 16).Value = 0.7927180626278817
$ws.Cells.Item(17, 17).Value = 1075.023482042533
$ws.Cells.Item(17, 18).Value = 9675.2113383828
$ws.Cells.Item(17, 19).Value = 0.01287229424337771
$ws.Cells.Item(17, 20).Value = 0.01287229424337771
